$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.146.59"
$ws.Range("E2").Value = "  +0.81%  "

$ws.Range("D3").Value = "2.570.22"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "505.80"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.60"
$ws.Range("E6").Value = "  -2.65%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("E8").Value = "  -4.86%  "

$ws.Range("D9").Value = "2.570.03"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.54"
$ws.Range("E10").Value = "  +7.34%  "

$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("E13").Value = "  +1.19%  "

$ws.Range("D14").Value = "3.019.07"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").Value = "60.152.28"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.50"
$ws.Range("E16").Value = "  -1.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("E17").Value = "  +1.93%  "

$ws.Range("D18").Value = "2.571.87"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.77"
$ws.Range("E19").Value = "  +0.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "344.98"
$ws.Range("E20").Value = "  +1.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.39"
$ws.Range("E21").Value = "  +0.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.09"
$ws.Range("E22").Value = "  +0.92%  "

$ws.Range("E23").Value = "  -0.41%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "59.68"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("E25").Value = "  +0.21%  "

$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").Value = "0.0₃0840"
$ws.Range("E28").Value = "  +0.96%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("E29").Value = "  +0.60%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.27"
$ws.Range("E31").Value = "  -0.79%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.10"
$ws.Range("E32").Value = "  -2.46%  "

$ws.Range("E33").Value = "  -0.90%  "

$ws.Range("E34").Value = "  +3.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.97"
$ws.Range("E35").Value = "  +1.84%  "

$ws.Range("E36").Value = "  -1.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.855"
$ws.Range("E37").Value = "  +9.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.848"
$ws.Range("E38").Value = "  -0.73%  "

$ws.Range("E39").Value = "  +2.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.04"
$ws.Range("E40").Value = "  +2.33%  "

$ws.Range("E41").Value = "  +0.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "292.18"
$ws.Range("E42").Value = "  -3.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.618"
$ws.Range("E43").Value = "  -1.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0993"
$ws.Range("E44").Value = "  -2.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.997"
$ws.Range("E45").Value = "  -0.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0554"
$ws.Range("E46").Value = "  -3.51%  "

$ws.Range("E47").Value = "  +2.32%  "

$ws.Range("E48").Value = "  -1.62%  "

$ws.Range("E49").Value = "  -1.75%  "

$ws.Range("E50").Value = "  +0.13%  "

$ws.Range("D51").Value = "1.992.80"
$ws.Range("E51").Value = "  +0.36%  "
